$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.703.32'
$ws.Range("E2").Value = '  -2.39%  '
$ws.Range("D3").Value = '1.871.67'
$ws.Range("E3").Value = '  -2.19%  '
$ws.Range("E4").Value = '  -0.89%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.02'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.09%  '
$ws.Range("E6").Value = '  -3.60%  '
$ws.Range("E7").Value = '  -0.91%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.24'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.346'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.19%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '50.18'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.69%  '
$ws.Range("E11").Value = '  +0.44%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0965'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.45%  '
$ws.Range("D13").Value = '2.141.51'
$ws.Range("E13").Value = '  -2.23%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '12.80'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.57%  '
$ws.Range("E15").Value = '  -0.91%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.87'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.92%  '
$ws.Range("D17").Value = '1.856.46'
$ws.Range("E17").Value = '  -3.16%  '
$ws.Range("D18").Value = '34.651.11'
$ws.Range("E18").Value = '  -2.44%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '72.70'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.05%  '
$ws.Range("E20").Value = '  -1.08%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '246.95'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.80%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.67'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.83%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.90'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.36%  '
$ws.Range("E24").Value = '  -0.92%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.39'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.79%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.19'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.97%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '164.62'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.59%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.34'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.79%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.18'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.18%  '
$ws.Range("E30").Value = '  -4.71%  '
$ws.Range("E31").Value = '  -0.14%  '
$ws.Range("E32").Value = '  -0.22%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0578'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.19%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.56'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.17%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.14'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.91%  '
$ws.Range("E36").Value = '  -0.94%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.82'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.92%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.829'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -9.78%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.98'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.76%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '17.19'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.70%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '97.64'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.35%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0659'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.57%  '
$ws.Range("E43").Value = '  -0.49%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.08'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.78%  '
$ws.Range("D45").Value = '1.288.46'
$ws.Range("E45").Value = '  -4.92%  '
$ws.Range("E46").Value = '  -4.99%  '
$ws.Range("E47").Value = '  -0.86%  '
$ws.Range("E48").Value = '  -2.41%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '12.13'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.33%  '
$ws.Range("E50").Value = '  +5.44%  '
$ws.Range("E51").Value = '  -2.04%  '
